$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.652441620826721
$ws.Range("B1").Value = 1.774739265441895
$ws.Range("C1").Value = 2.013420581817627
$ws.Range("D1").Value = 3.162654399871826
$ws.Range("E1").Value = 3.542757511138916
